$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# wdFindStop = 0 (do not wrap outside the supplied range)
# wdReplaceOne = 1 (replace just the single match found)

# --- Row 6 = Version Control row "2.8.1" ---
# "ผู้รับผิดชอบ" (responsible) cell: ณัฐนันท์ (QA)  ->  ณัฐดนัย (DM)
$rng = $t.Rows.Item(6).Cells.Item(4).Range
$res = $rng.Find.Execute("ณัฐนันท์", $true, $false, $false, $false, $false, $true, 0, $false, "ณัฐดนัย", 1)
$rng = $t.Rows.Item(6).Cells.Item(4).Range
$res = $rng.Find.Execute(" (QA)", $true, $false, $false, $false, $false, $true, 0, $false, " (DM)", 1)

# "ผู้ตรวจ" (reviewer) cell: กิตติพศ (SP)  ->  วิรัตน์ (TL)
$rng = $t.Rows.Item(6).Cells.Item(5).Range
$res = $rng.Find.Execute("กิตติพศ ", $true, $false, $false, $false, $false, $true, 0, $false, "วิรัตน์", 1)
$rng = $t.Rows.Item(6).Cells.Item(5).Range
$res = $rng.Find.Execute("(SP)", $true, $false, $false, $false, $false, $true, 0, $false, " (TL)", 1)

# --- Row 7 = Version Control row "2.4.1" ---
# "ผู้รับผิดชอบ" (responsible) cell: ณัฐนันท์ (QA)  ->  วิรัตน์ (TL)
$rng = $t.Rows.Item(7).Cells.Item(4).Range
$res = $rng.Find.Execute("ณัฐนันท์", $true, $false, $false, $false, $false, $true, 0, $false, "วิรัตน์", 1)
$rng = $t.Rows.Item(7).Cells.Item(4).Range
$res = $rng.Find.Execute(" (QA)", $true, $false, $false, $false, $false, $true, 0, $false, " (TL)", 1)

# "ผู้ตรวจ" (reviewer) cell: กิตติพศ (SP)  ->  วริศรา (D)
$rng = $t.Rows.Item(7).Cells.Item(5).Range
$res = $rng.Find.Execute("กิตติพศ ", $true, $false, $false, $false, $false, $true, 0, $false, "วริศรา", 1)
$rng = $t.Rows.Item(7).Cells.Item(5).Range
$res = $rng.Find.Execute("(SP)", $true, $false, $false, $false, $false, $true, 0, $false, " (D)", 1)

Write-Output "done"
